# Applies the "Updated cryptos list" refresh: new Price (col D) and
# Volume(1h) (col E) figures for each coin row, plus three coins whose
# row order changed (Solana/WrappedEther and RenderToken/Decentraland/
# NEARProtocol), which also updates their Coin (B) and Link (C) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Writes $text as a literal string even when it looks like a number
    # (e.g. "1.022" or "82.00"): a leading apostrophe forces Excel to
    # store it as text instead of auto-converting/rounding it, and
    # resetting the style back to Normal avoids leaving a stray
    # quote-prefixed number format behind on the cell.
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.345.89"
$ws.Range("E2").Value = "  +3.80%  "
# Row 3
$ws.Range("D3").Value = "1.838.21"
$ws.Range("E3").Value = "  +4.00%  "
# Row 4
$ws.Range("E4").Value = "  +3.01%  "
# Row 5
Set-TextValue $ws.Range("D5") "319.61"
$ws.Range("E5").Value = "  +4.65%  "
# Row 6
Set-TextValue $ws.Range("D6") "1.022"
$ws.Range("E6").Value = "  +2.63%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.4346"
$ws.Range("E7").Value = "  +1.69%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.3717"
$ws.Range("E8").Value = "  +2.20%  "
# Row 9
Set-TextValue $ws.Range("D9") "0.07326"
$ws.Range("E9").Value = "  +2.18%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.8749"
$ws.Range("E10").Value = "  +3.16%  "
# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D11") "21.34"
$ws.Range("E11").Value = "  +4.95%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.032.39"
$ws.Range("E12").Value = "  +14.60%  "
# Row 13
$ws.Range("E13").Value = "  +4.34%  "
# Row 14
Set-TextValue $ws.Range("D14") "6.675"
$ws.Range("E14").Value = "  +3.90%  "
# Row 15
Set-TextValue $ws.Range("D15") "0.07155"
$ws.Range("E15").Value = "  +4.59%  "
# Row 16
Set-TextValue $ws.Range("D16") "82.00"
$ws.Range("E16").Value = "  +4.16%  "
# Row 17
Set-TextValue $ws.Range("D17") "1.028"
$ws.Range("E17").Value = "  +2.73%  "
# Row 18
Set-TextValue $ws.Range("D18") "0.000009003"
$ws.Range("E18").Value = "  +3.87%  "
# Row 19
Set-TextValue $ws.Range("D19") "1.019"
$ws.Range("E19").Value = "  +2.21%  "
# Row 20
$ws.Range("E20").Value = "  +3.06%  "
# Row 21
$ws.Range("D21").Value = "27.384.52"
$ws.Range("E21").Value = "  +3.91%  "
# Row 22
Set-TextValue $ws.Range("D22") "5.238"
$ws.Range("E22").Value = "  +2.88%  "
# Row 23
$ws.Range("E23").Value = "  +0.08%  "
# Row 24
$ws.Range("D24").Value = "2.225.99"
$ws.Range("E24").Value = "  +11.87%  "
# Row 25
Set-TextValue $ws.Range("D25") "156.55"
$ws.Range("E25").Value = "  +3.33%  "
# Row 26
Set-TextValue $ws.Range("D26") "1.902"
$ws.Range("E26").Value = "  +2.53%  "
# Row 27
$ws.Range("E27").Value = "  +2.69%  "
# Row 28
Set-TextValue $ws.Range("D28") "5.279"
$ws.Range("E28").Value = "  +3.70%  "
# Row 29
Set-TextValue $ws.Range("D29") "1.930"
$ws.Range("E29").Value = "  +7.10%  "
# Row 30
Set-TextValue $ws.Range("D30") "115.37"
$ws.Range("E30").Value = "  +1.53%  "
# Row 31
Set-TextValue $ws.Range("D31") "0.09000"
$ws.Range("E31").Value = "  +0.78%  "
# Row 32
$ws.Range("E32").Value = "  +6.58%  "
# Row 33
Set-TextValue $ws.Range("D33") "0.7583"
$ws.Range("E33").Value = "  +4.10%  "
# Row 34
Set-TextValue $ws.Range("D34") "4.458"
$ws.Range("E34").Value = "  +3.31%  "
# Row 35
Set-TextValue $ws.Range("D35") "2.848"
$ws.Range("E35").Value = "  +4.76%  "
# Row 36
Set-TextValue $ws.Range("D36") "1.024"
$ws.Range("E36").Value = "  +2.89%  "
# Row 37
Set-TextValue $ws.Range("D37") "1.149"
$ws.Range("E37").Value = "  +5.38%  "
# Row 38
Set-TextValue $ws.Range("D38") "0.01951"
$ws.Range("E38").Value = "  +3.15%  "
# Row 39
Set-TextValue $ws.Range("D39") "0.05260"
$ws.Range("E39").Value = "  +2.14%  "
# Row 40
Set-TextValue $ws.Range("D40") "0.5161"
$ws.Range("E40").Value = "  +4.88%  "
# Row 41
Set-TextValue $ws.Range("D41") "2.798"
$ws.Range("E41").Value = "  +8.50%  "
# Row 42
$ws.Range("E42").Value = "  +3.23%  "
# Row 43
Set-TextValue $ws.Range("D43") "6.514"
$ws.Range("E43").Value = "  +3.72%  "
# Row 44
Set-TextValue $ws.Range("D44") "8.450"
$ws.Range("E44").Value = "  +5.45%  "
# Row 45
Set-TextValue $ws.Range("D45") "107.85"
$ws.Range("E45").Value = "  +2.95%  "
# Row 46
Set-TextValue $ws.Range("D46") "10.53"
$ws.Range("E46").Value = "  +3.80%  "
# Row 47
Set-TextValue $ws.Range("D47") "1.025"
$ws.Range("E47").Value = "  +3.03%  "
# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.4620"
$ws.Range("E48").Value = "  +3.01%  "
# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.664"
$ws.Range("E49").Value = "  +3.32%  "
# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "1.898"
$ws.Range("E50").Value = "  +9.25%  "
# Row 51
Set-TextValue $ws.Range("D51") "0.06282"
$ws.Range("E51").Value = "  +1.66%  "
